$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4's audio record was mislabeled "Лекция 2"; correct it to "Лекция 1".
# This removes the now-unused "Лекция 2" entry from the shared strings table.
$ws.Range("A4").Value = "Лекция 1"

# Move/update the current selection from D6 to A4.
$ws.Range("A4").Select()
